# Generate Report for Handoff
#
# The localization-status report was regenerated: a new handoff run replaced
# the old source-file UUID (e602c46e-7f9a-43aa-8dac-10c72f177e62) with a new
# one (3d9b462a-26a7-468f-ab45-1a36075eaa60), and produced new handoff
# artifact names + timestamps for the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newId = "3d9b462a-26a7-468f-ab45-1a36075eaa60"

# Source file name (A2) is the same on all three sheets.
$overview.Range("A2").Value2 = "$newId.md"
$zhcn.Range("A2").Value2     = "$newId.md"
$dede.Range("A2").Value2     = "$newId.md"

# zh-cn: new handoff file name + new handoff datetime.
$zhcn.Range("C2").Value2 = "$newId.ec8657d41dbe20158d00d54696241c97662d68cc.zh-cn.xlf"
$zhcn.Range("D2").Value2 = "2016-03-09 08:38:35"

# de-de: new handoff file name + new handoff datetime.
$dede.Range("C2").Value2 = "$newId.ec8657d41dbe20158d00d54696241c97662d68cc.de-de.xlf"
$dede.Range("D2").Value2 = "2016-03-09 08:38:40"
